$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price (D) and volume-change (E) columns per the GitHub Actions refresh.
# D-column values must stay plain text (prices are formatted like "27.362.91" or "1.0000" with
# significant trailing zeros), so we force a Text number format while writing, then restore the
# original (unstyled) look so no stray style index is left behind.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.362.91'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.84%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.825.48'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.93%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.0000'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.25%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.37'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.02%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9996'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.24%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4476'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.06%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3771'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.03%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07405'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.01%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8794'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.40%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.91'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.75%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.825.50'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.00%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.719'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.94%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.431'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.71%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '92.68'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.44%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.07054'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.54%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.001'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.29%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008805'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.99%  '

$ws.Range("E19").Value = '  -0.25%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.05'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.43%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '27.374.74'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.78%  '

$ws.Range("E22").Value = '  +4.36%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.95'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.84%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.961'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.05%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '150.80'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.54%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.284'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.71%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.61'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.40%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.350'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.10%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '117.24'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.09%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.08876'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.64%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.7902'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +6.08%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.197'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.23%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.571'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.51%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.928'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.32%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9991'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.28%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.110'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.28%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01977'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.41%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05271'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.68%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '7.314'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.30%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5307'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.80%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.880'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.66%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.323'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +19.22%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1702'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.03%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.639'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.92%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5052'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.67%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.56'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.11%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '105.44'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.31%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.686'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.14%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.9989'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.23%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06401'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.75%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '66.58'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +6.63%  '
